# Added usecases for multithreading.
# Append two new candidate rows (sangeeta sahu / Ravi kant sharma) to the
# "JavaJMS" use-case tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: S.No=2, Name=sangeeta sahu, Committed=Yes, Commit Date&Time
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "sangeeta sahu"
$ws.Range("C5").Value = "Yes"
$ws.Range("D5").Value = 42711.684027777781
$ws.Range("D5").NumberFormat = "m/d/yy h:mm"

# New row 6: S.No=3, Name=Ravi kant sharma, Committed=Yes, Commit Date&Time
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Ravi kant sharma"
$ws.Range("C6").Value = "Yes"
$ws.Range("D6").Value = 42711.684027777781
$ws.Range("D6").NumberFormat = "m/d/yy h:mm"

# Match the author's final selection/cursor position recorded in the sheet.
$ws.Range("D9").Select()
